$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in rows 218-222 ---
$ws.Range("F218").Value = 110
$ws.Range("L218").Value = -4286
$ws.Range("M218").Value = -142
$ws.Range("Q218").Value = 2015
$ws.Range("S218").Value = 2294
$ws.Range("U218").Value = -152
$ws.Range("V218").Value = 2506
$ws.Range("W218").Value = 1673
$ws.Range("X218").Value = 906
$ws.Range("Z218").Value = -72
$ws.Range("AA218").Value = -669
$ws.Range("AB218").Value = -1210
$ws.Range("AC218").Value = 96
$ws.Range("AE218").Value = 448
$ws.Range("F219").Value = 244
$ws.Range("L219").Value = -241
$ws.Range("M219").Value = 44
$ws.Range("Q219").Value = 105
$ws.Range("S219").Value = -133
$ws.Range("U219").Value = 343
$ws.Range("V219").Value = 61
$ws.Range("W219").Value = -177
$ws.Range("X219").Value = 307
$ws.Range("Z219").Value = -68
$ws.Range("AA219").Value = 1111
$ws.Range("AB219").Value = 144
$ws.Range("AC219").Value = 837
$ws.Range("AE219").Value = 115
$ws.Range("B220").Value = -734
$ws.Range("E220").Value = -906
$ws.Range("F220").Value = 1681
$ws.Range("M220").Value = -318
$ws.Range("Q220").Value = -1161
$ws.Range("S220").Value = -918
$ws.Range("U220").Value = -196
$ws.Range("V220").Value = 796
$ws.Range("W220").Value = 400
$ws.Range("X220").Value = 481
$ws.Range("Z220").Value = -89
$ws.Range("AA220").Value = -1538
$ws.Range("AB220").Value = -921
$ws.Range("AC220").Value = -219
$ws.Range("AE220").Value = -466
$ws.Range("F221").Value = 2998
$ws.Range("Q221").Value = 47
$ws.Range("S221").Value = 356
$ws.Range("V221").Value = 458
$ws.Range("W221").Value = 284
$ws.Range("X221").Value = 200
$ws.Range("Z221").Value = -25
$ws.Range("AA221").Value = -1201
$ws.Range("AB221").Value = -3397
$ws.Range("AC221").Value = -667
$ws.Range("AE221").Value = 2844
$ws.Range("F222").Value = -1278
$ws.Range("L222").Value = 1003
$ws.Range("M222").Value = -57
$ws.Range("P222").Value = 514
$ws.Range("Q222").Value = -7159
$ws.Range("S222").Value = -7385
$ws.Range("V222").Value = -153
$ws.Range("W222").Value = -285
$ws.Range("X222").Value = 75
$ws.Range("Z222").Value = 71
$ws.Range("AA222").Value = 1384
$ws.Range("AB222").Value = 825
$ws.Range("AC222").Value = 751
$ws.Range("AE222").Value = -229

# --- Add new row 223 ---
# A223 needs to be the text "01-06-2021" (not auto-converted to a date).
# Enter it as a formula producing the text, then paste-special as values
# so the stored cell is a plain shared string with no style changes.
$ws.Cells.Item(223, 1).Formula = '="01-06-2021"'
$ws.Cells.Item(223, 1).Copy()
$ws.Cells.Item(223, 1).PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("B223").Value = -4090
$ws.Range("C223").Value = -3365
$ws.Range("D223").Value = 0
$ws.Range("E223").Value = -725
$ws.Range("F223").Value = 893
$ws.Range("G223").Value = -2608
$ws.Range("H223").Value = -2596
$ws.Range("I223").Value = -3
$ws.Range("J223").Value = 0
$ws.Range("K223").Value = -9
$ws.Range("L223").Value = 1782
$ws.Range("M223").Value = -188
$ws.Range("N223").Value = 208
$ws.Range("O223").Value = 0
$ws.Range("P223").Value = 1762
$ws.Range("Q223").Value = 490
$ws.Range("R223").Value = 0
$ws.Range("S223").Value = 419
$ws.Range("T223").Value = 42
$ws.Range("U223").Value = 29
$ws.Range("V223").Value = 1229
$ws.Range("W223").Value = 62
$ws.Range("X223").Value = 1200
$ws.Range("Y223").Value = -9
$ws.Range("Z223").Value = -25
$ws.Range("AA223").Value = 904
$ws.Range("AB223").Value = -676
$ws.Range("AC223").Value = 898
$ws.Range("AD223").Value = 59
$ws.Range("AE223").Value = 623
